# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos table
# with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '62.560.95'
    'E2' = '  +4.42%  '
    'D3' = '2.427.59'
    'E3' = '  +5.22%  '
    'E4' = '  -0.02%  '
    'D5' = '558.12'
    'E5' = '  +3.06%  '
    'D6' = '139.15'
    'E6' = '  +7.40%  '
    'E7' = '  -0.04%  '
    'E8' = '  +2.09%  '
    'D9' = '2.426.10'
    'E9' = '  +5.23%  '
    'E10' = '  +3.76%  '
    'E11' = '  +4.42%  '
    'E12' = '  +0.43%  '
    'E13' = '  +4.80%  '
    'D14' = '26.21'
    'E14' = '  +12.36%  '
    'D15' = '2.860.66'
    'E15' = '  +5.23%  '
    'D16' = '62.428.74'
    'E16' = '  +4.21%  '
    'E17' = '  +7.33%  '
    'D18' = '2.433.09'
    'E18' = '  +5.69%  '
    'D19' = '11.24'
    'E19' = '  +7.21%  '
    'D20' = '346.11'
    'E20' = '  +10.73%  '
    'D21' = '4.20'
    'E21' = '  +3.00%  '
    'D22' = '6.84'
    'E22' = '  +4.26%  '
    'E23' = '  -0.04%  '
    'E24' = '  -2.81%  '
    'D25' = '65.34'
    'E25' = '  +2.67%  '
    'E26' = '  +1.53%  '
    'D27' = '1.57'
    'E27' = '  +16.00%  '
    'E28' = '  -0.03%  '
    'E29' = '  +5.38%  '
    'D30' = '1.36'
    'E30' = '  +15.83%  '
    'E31' = '  +5.61%  '
    'E32' = '  +8.17%  '
    'E33' = '  +10.97%  '
    'D34' = '172.62'
    'E34' = '  +0.66%  '
    'E35' = '  +6.24%  '
    'E36' = '  +4.90%  '
    'D37' = '378.78'
    'E37' = '  +19.79%  '
    'D38' = '18.58'
    'E38' = '  +5.06%  '
    'D39' = '4.46'
    'E39' = '  +11.48%  '
    'D41' = '0.999'
    'E41' = '  -0.14%  '
    'E42' = '  +12.62%  '
    'D43' = '39.44'
    'E43' = '  +3.54%  '
    'D44' = '144.82'
    'E44' = '  +6.44%  '
    'D45' = '3.68'
    'E45' = '  +7.41%  '
    'D46' = '20.85'
    'E46' = '  +10.93%  '
    'D48' = '0.0953'
    'E48' = '  +1.62%  '
    'D49' = '0.0521'
    'E49' = '  +6.14%  '
    'E50' = '  +4.62%  '
    'D51' = '17.84'
    'E51' = '  +6.42%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    if ($addr[0] -eq "D") {
        # Price column holds numeric-looking text (e.g. "62.560.95",
        # "0.999"); force text storage so Excel does not reinterpret it
        # as a number/date, then drop back to the default style so no
        # visible formatting change is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$addr]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$addr]
    }
}
